# Lightning Talk (Collections) - "Code added" update
#
# Slide 11 ("Thank you for the attention") has a "Code:" line and a
# "Slides:" line, each built from 4 runs:
#   "https://github.com/greenfox-academy/" + "bramble100" (bold)
#   + "/tree/master/week-03/" + "<repo-name>" (bold)
#
# Before the edit the "Code:" line's hyperlink/text still pointed at last
# week's "SortablePokerHands" repo, while the "Slides:" line already
# pointed at this week's "lightning-talk" repo (a different relationship
# id). The commit retargets the "Code:" line's last run text to
# "lightning-talk" and repoints the "Slides:" line's hyperlinks onto the
# same relationship ("rId2") that the "Code:" line already uses.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- "Code:" line (paragraph 6): last run text SortablePokerHands -> lightning-talk
$codeLinkPara = $tr.Paragraphs(6, 1)
$codeLastRun = $codeLinkPara.Runs(4, 1)
$codeLastRun.Text = "lightning-talk"

# --- "Slides:" line (paragraph 8): repoint all 4 runs' hyperlink onto the
#     same target the "Code:" line uses, so they end up sharing its
#     relationship id instead of keeping their own.
$slidesLinkPara = $tr.Paragraphs(8, 1)
$targetAddress = "https://github.com/greenfox-academy/bramble100/tree/master/week-02/SortablePokerHands"

for ($i = 1; $i -le 4; $i++) {
  $run = $slidesLinkPara.Runs($i, 1)
  $run.ActionSettings(1).Hyperlink.Address = ""
}
for ($i = 1; $i -le 4; $i++) {
  $run = $slidesLinkPara.Runs($i, 1)
  $run.ActionSettings(1).Hyperlink.Address = $targetAddress
}
